$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = 484513.94
$ws.Range("D5").Value = 128787.19
$ws.Range("D4").Select()
